$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Des2"
$ws.Range("B3").Value = "Car"
$ws.Range("C3").Value = 2021
$ws.Range("D3").Value = "Mercedes"
$ws.Range("E3").Value = "A6"
$ws.Range("F3").Value = 982514
$ws.Range("G3").Value = 258
$ws.Range("H3").Value = "Gasoline"
$ws.Range("I3").Value = "No"
$ws.Range("J3").Value = 6

$ws.Range("A4").Value = "Des 3"
$ws.Range("B4").Value = "Car"
$ws.Range("C4").Value = 2022
$ws.Range("D4").Value = "Audi"
$ws.Range("E4").Value = "Q2"
$ws.Range("F4").Value = 98512
$ws.Range("G4").Value = 258
$ws.Range("H4").Value = "Gasoline"
$ws.Range("I4").Value = "Yes"
$ws.Range("J4").Value = 5

$ws.Range("F4").Select()
